# Delete row 15 ("Physics" / "Implement a separate collisionbox that
# autocrops on load") from the "Todo" sheet. Excel automatically shifts
# the remaining rows up, updates the AutoFilter range and the
# _xlnm._FilterDatabase defined name, and removes the now-unused shared
# string on save.

$wb = $excel.ActiveWorkbook
$todo = $wb.Worksheets.Item("Todo")
$todo.Activate()

$todo.Rows.Item(15).Delete()

# Match the recorded selection/scroll position after the edit.
$todo.Application.ActiveWindow.ScrollRow = 17
$todo.Range("C29").Select()
